$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9043993353843689
$ws.Range("B1").Value = 4.112283706665039
$ws.Range("C1").Value = 6.101568698883057
$ws.Range("D1").Value = 1.494413733482361
$ws.Range("E1").Value = 0.8417622447013855
